# "Running a2,a4 and b6" — mark all B-suite test cases back to "not yet
# run" (Runmode = N) except TestCase_B6 (row 7), which has now been run
# and passed (Results = PASS). Selection moves to the next block to run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Reset Runmode to "N" for every row except row 7 (TestCase_B6).
$ws.Range("C2:C6").Value = "N"
$ws.Range("C8:C28").Value = "N"

# Row 7 (TestCase_B6) has been executed: mark its Results as PASS.
$ws.Range("D7").Value = "PASS"

# Reflect the in-progress selection left by the run (next rows to run).
$ws.Range("C8:C28").Select()
